# Fills in the "test" placeholder runs across every certificate slide with
# the real values: recipient name "REC Rick", signer name "Jazz", and the
# effective date "30th Aug 2023".
#
# Each slide repeats the same 7 placeholder occurrences (in shape/document
# order):
#   1) bare "test"            -> "REC Rick"   (recipient name, big italic)
#   2) "MAJ test"/"LTCOL test" -> "... Jazz"   (signer rank + name)
#   3) same as #2 in a different font run      -> "... Jazz"
#   4) "test,"/"Major test," etc.              -> "... Jazz,"  (signer name)
#   5) bare "test"            -> "REC Rick"   (recipient name again)
#   6) "test" / "test."       -> "30th Aug 2023" (+ trailing punctuation)
#   7) "test" / "test "       -> "REC Rick"   (+ trailing space)

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $occurrence = 0

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        if (-not $shp.TextFrame.HasText) { continue }

        $tr = $shp.TextFrame.TextRange
        $whole = $tr.Text
        if ($whole.IndexOf("test") -lt 0) { continue }

        $occurrence = $occurrence + 1

        $idx = $whole.IndexOf("test")
        while ($idx -ge 0) {
            if ($occurrence -eq 6) {
                $replacement = "30th Aug 2023"
            } elseif ($occurrence -eq 2 -or $occurrence -eq 3 -or $occurrence -eq 4) {
                $replacement = "Jazz"
            } else {
                $replacement = "REC Rick"
            }

            $start = $idx + 1
            $len = 4
            $sub = $tr.Characters($start, $len)
            $sub.Text = $replacement

            $whole = $tr.Text
            $searchFrom = $idx + $replacement.Length
            if ($searchFrom -ge $whole.Length) {
                $idx = -1
            } else {
                $idx = $whole.IndexOf("test", $searchFrom)
            }
        }
    }
}
